$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.206.77'
$ws.Range('E2').Value = '  -1.67%  '
$ws.Range('D3').Value = '3.118.65'
$ws.Range('E3').Value = '  -2.55%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '595.18'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('D6').Value = '157.99'
$ws.Range('E6').Value = '  +2.66%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = '3.119.40'
$ws.Range('E9').Value = '  -2.42%  '
$ws.Range('E10').Value = '  -4.86%  '
$ws.Range('D11').Value = '5.91'
$ws.Range('E11').Value = '  -3.19%  '
$ws.Range('E12').Value = '  -3.93%  '
$ws.Range('D13').Value = '37.23'
$ws.Range('E13').Value = '  -5.51%  '
$ws.Range('E14').Value = '  -5.54%  '
$ws.Range('D15').Value = '3.634.37'
$ws.Range('E15').Value = '  -2.62%  '
$ws.Range('E16').Value = '  -1.54%  '
$ws.Range('D17').Value = '7.25'
$ws.Range('E17').Value = '  -2.25%  '
$ws.Range('D18').Value = '64.132.82'
$ws.Range('D19').Value = '3.120.04'
$ws.Range('E19').Value = '  -2.50%  '
$ws.Range('D20').Value = '477.47'
$ws.Range('E20').Value = '  -1.35%  '
$ws.Range('D21').Value = '14.53'
$ws.Range('E21').Value = '  -3.85%  '
$ws.Range('D22').Value = '0.716'
$ws.Range('E22').Value = '  -7.44%  '
$ws.Range('E23').Value = '  -4.41%  '
$ws.Range('D24').Value = '2.49'
$ws.Range('E24').Value = '  +2.00%  '
$ws.Range('E25').Value = '  -6.88%  '
$ws.Range('D26').Value = '81.40'
$ws.Range('E26').Value = '  -2.69%  '
$ws.Range('D27').Value = '10.60'
$ws.Range('E27').Value = '  +7.50%  '
$ws.Range('E28').Value = '  -0.30%  '
$ws.Range('D29').Value = '7.61'
$ws.Range('E29').Value = '  +1.44%  '
$ws.Range('E30').Value = '  -2.59%  '
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('D32').Value = '2.20'
$ws.Range('E32').Value = '  -3.54%  '
$ws.Range('E33').Value = '  -6.25%  '
$ws.Range('D34').Value = '27.34'
$ws.Range('E34').Value = '  -4.35%  '
$ws.Range('E35').Value = '  -5.41%  '
$ws.Range('D36').Value = '1.07'
$ws.Range('E36').Value = '  -2.33%  '
$ws.Range('E37').Value = '  -4.98%  '
$ws.Range('E38').Value = '  -7.12%  '
$ws.Range('E39').Value = '  -5.37%  '
$ws.Range('D40').Value = '51.00'
$ws.Range('E40').Value = '  -1.06%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').Value = '449.50'
$ws.Range('E41').Value = '  -5.33%  '
$ws.Range('B42').Value = 'Cosmos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D42').Value = '9.16'
$ws.Range('E42').Value = '  -3.34%  '
$ws.Range('D43').Value = '0.294'
$ws.Range('E43').Value = '  -2.78%  '
$ws.Range('D44').Value = '0.0367'
$ws.Range('E44').Value = '  -4.34%  '
$ws.Range('B45').Value = 'Arweave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D45').Value = '40.65'
$ws.Range('E45').Value = '  +5.23%  '
$ws.Range('B46').Value = 'Kaspa'
$ws.Range('C46').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D46').Value = '0.112'
$ws.Range('E46').Value = '  +0.37%  '
$ws.Range('D47').Value = '2.832.68'
$ws.Range('E47').Value = '  -4.37%  '
$ws.Range('D48').Value = '130.87'
$ws.Range('E48').Value = '  -0.65%  '
$ws.Range('D49').Value = '25.93'
$ws.Range('E49').Value = '  +0.87%  '
$ws.Range('E51').Value = '  -3.11%  '
